$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- EDU_LEVEL (row 5) algorithm update ---
# Editing H5 in place with the revised case_when() text. Because this string
# is uniquely referenced, the shared-string table slot it occupied is
# dropped and the new text is appended at the end, which is what shifts the
# neighbouring __BLANK__ / paste / label strings down by one index (and is
# exactly what the target workbook shows).
$newFormula = @'
case_when(
  m_berufab == 5 | v_berufab == 5 ~ 7,
  m_berufab %in% c(3, 4) | v_berufab %in% c(3, 4) ~ 6,
  m_berufab %in% c(1, 2) | v_berufab %in% c(1, 2) ~ 4,
  m_schulab %in% c(3, 4) | v_schulab %in% c(3, 4) ~ 3,
  m_schulab %in% c(1, 2) | v_schulab %in% c(1, 2) ~ 2,
  m_schulab == 5 | v_schulab == 5 | m_berufab %in% c(6, 7, 8)| v_berufab %in% c(6, 7, 8) ~ 9,
  TRUE ~ NA_real_
)
'@
$ws.Range("H5").Value = $newFormula

# Row 5 wraps the (now longer) algorithm text; Excel auto-fit pins it at the
# worksheet's maximum row height.
$ws.Rows.Item(5).RowHeight = 409.5

# --- status_detail ("identical" -> "compatible") for the P2 food-group rows ---
$foodGroupRows = @(97, 98, 99, 100, 102, 103, 104, 105, 106, 107, 108)
foreach ($r in $foodGroupRows) {
    $ws.Cells.Item($r, 11).Value = "compatible"
}

# --- DIETARY_ASSESS_INSTR (row 109): algorithm value 1 -> 3, drop red font ---
$ws.Cells.Item(109, 8).Value = 3
$ws.Range("H101").Copy()
$ws.Range("H109").PasteSpecial(-4122)   # xlPasteFormats

# --- Selection / scroll position ---
$ws.Range("F108").Select()
